$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are purely numeric-looking strings need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (the source workbook stores every Price/Volume cell as text).
$ws.Range('D2').Value = '63.565.42'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.41'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.77'
$ws.Range('E6').Value = '  +3.52%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +6.56%  '
$ws.Range('D9').Value = '3.079.14'
$ws.Range('E9').Value = '  -2.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.156'
$ws.Range('E10').Value = '  -3.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.85'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000243'
$ws.Range('E13').Value = '  -2.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.52'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '3.591.85'
$ws.Range('E16').Value = '  -2.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.20'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').Value = '63.558.53'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '3.082.85'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.95'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.720'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.56'
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.14'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.04'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.44'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = '0.0₃0854'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.05'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.37'
$ws.Range('E37').Value = '  +2.78%  '
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.22'
$ws.Range('E39').Value = '  -4.84%  '
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.60'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '445.17'
$ws.Range('E42').Value = '  -5.17%  '
$ws.Range('E43').Value = '  -5.23%  '
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('D45').Value = '2.815.74'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.52'
$ws.Range('E47').Value = '  -3.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.43'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.113'
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.15'
$ws.Range('E51').Value = '  +2.54%  '
